# Admin - Manage Test - Pretest and Posttest - Edit
# Adds two new worksheets at the end of the workbook:
#   "Admin-MT-PrePost- Delete"  (search/delete pretest-posttest questions)
#   "Admin-MT-PrePost- Edit"    (search/edit pretest-posttest questions)

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- New sheet: Admin-MT-PrePost- Delete -----------------------------------
$wsDelete = $wb.Worksheets.Add($null, $lastSheet)
$wsDelete.Name = "Admin-MT-PrePost- Delete"

$wsDelete.Columns.Item(3).ColumnWidth = 10.0

$wsDelete.Range("A1").Value = "searchID"
$wsDelete.Range("B1").Value = "kondisi"
$wsDelete.Range("C1").Value = "keterangan"

$wsDelete.Range("A2").Value = "qwerty"
$wsDelete.Range("B2").Value = "fail"
$wsDelete.Range("C2").Value = "notFound"

$wsDelete.Range("A3").Value = 1771
$wsDelete.Range("A3").HorizontalAlignment = -4131
$wsDelete.Range("B3").Value = "pass"
$wsDelete.Range("C3").Value = "cancel"

$wsDelete.Range("A4").Value = 1771
$wsDelete.Range("A4").HorizontalAlignment = -4131
$wsDelete.Range("B4").Value = "pass"

$wsDelete.Range("A5").Select()

# --- New sheet: Admin-MT-PrePost- Edit --------------------------------------
$wsEdit = $wb.Worksheets.Add($null, $wsDelete)
$wsEdit.Name = "Admin-MT-PrePost- Edit"

$wsEdit.Columns.Item(4).ColumnWidth = 10.0
$wsEdit.Columns.Item(5).ColumnWidth = 9.8333
$wsEdit.Columns.Item(7).ColumnWidth = 18.0
$wsEdit.Columns.Item(8).ColumnWidth = 14.8333
$wsEdit.Columns.Item(10).ColumnWidth = 12.6667
$wsEdit.Columns.Item(11).ColumnWidth = 17.6667
$wsEdit.Columns.Item(14).ColumnWidth = 10.0
$wsEdit.Columns.Item(15).ColumnWidth = 9.6667
$wsEdit.Columns.Item(16).ColumnWidth = 12.1667
$wsEdit.Columns.Item(18).ColumnWidth = 10.8333

$wsEdit.Range("A1").Value = "serachID"
$wsEdit.Range("B1").Value = "category"
$wsEdit.Range("C1").Value = "tipeSoal"
$wsEdit.Range("D1").Value = "jobFunction"
$wsEdit.Range("E1").Value = "jobPosition"
$wsEdit.Range("F1").Value = "module"
$wsEdit.Range("G1").Value = "technicalCompetence"
$wsEdit.Range("H1").Value = "levelCompetence"
$wsEdit.Range("I1").Value = "soal"
$wsEdit.Range("J1").Value = "jawabanUraian"
$wsEdit.Range("K1").Value = "jawabanBenarSalah"
$wsEdit.Range("L1").Value = "jawabanA"
$wsEdit.Range("M1").Value = "jawabanB"
$wsEdit.Range("N1").Value = "jawabanC"
$wsEdit.Range("O1").Value = "jawabanD"
$wsEdit.Range("P1").Value = "jawabanABCD"
$wsEdit.Range("Q1").Value = "kondisi"
$wsEdit.Range("R1").Value = "keterangan"

$wsEdit.Activate()
